$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add three new rows at the bottom of the transaction log (rows 73-75):
#   row 73 - blank separator (same look as the other separator rows in the
#            sheet, e.g. row 70), just carries the grey "date" formatting in
#            columns B:C
#   row 74 - a Payment transaction
#   row 75 - a Sale transaction (new payee -> needs a new shared string)
# The formatting used throughout the bottom of the sheet (grey font, date
# number format) lives in existing style indexes, so we copy the format from
# row 72 (the row right above the new block) instead of re-creating it, which
# keeps styles.xml untouched - exactly like Excel does when a user fills a
# new row down from the row above.
# ---------------------------------------------------------------------------

# Row 73: blank separator row - only B73:C73 carry the (date) style, no values
$ws.Range("B72:C72").Copy($ws.Range("B73:C73"))
$ws.Range("B73:C73").ClearContents()

# Row 74: Payment
$ws.Range("A72:E72").Copy($ws.Range("A74:E74"))
$ws.Range("A74").Value = "Payment"
$ws.Range("B74").Value = 42953
$ws.Range("C74").Value = 42953
$ws.Range("D74").Value = "AUTOMATIC PAYMENT - THANK"
$ws.Range("E74").Value = 596.68

# Row 75: Sale
$ws.Range("A72:E72").Copy($ws.Range("A75:E75"))
$ws.Range("A75").Value = "Sale"
$ws.Range("B75").Value = 42960
$ws.Range("C75").Value = 42961
$ws.Range("D75").Value = "TRADER JOE'S #196  QPS"
$ws.Range("E75").Value = -25.37

# ---------------------------------------------------------------------------
# Refresh the sort state over the newly added block (mirrors selecting
# A74:E77 and running Data > Sort by Post Date, which is how this sheet's
# sortState metadata gets maintained as rows are appended).
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B74:B77"))
$ws.Sort.SetRange($ws.Range("A74:E77"))
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# Update the view so the window is scrolled to show the new rows, with the
# cursor left on H76 (one row below the last entry), matching the saved
# worksheet view.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 66
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H76").Select()
